$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.245.24"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "3.202.56"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.14"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.96"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "3.203.17"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.66"
$ws.Range("E11").Value = "  -3.96%  "
$ws.Range("E12").Value = "  -3.44%  "
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "3.727.38"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "66.380.27"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.27"
$ws.Range("E17").Value = "  -3.34%  "
$ws.Range("D18").Value = "3.198.12"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "505.65"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.61"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.07"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.00"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.02"
$ws.Range("E28").Value = "  -2.78%  "
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  +42.70%  "
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "28.18"
$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.42"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.35"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "499.94"
$ws.Range("E38").Value = "  -1.99%  "
$ws.Range("D39").Value = "0.0₃0769"
$ws.Range("E39").Value = "  +12.31%  "
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("E41").Value = "  +5.25%  "
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("E44").Value = "  -2.12%  "
$ws.Range("D45").Value = "2.909.35"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.13"
$ws.Range("E47").Value = "  -1.64%  "
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.15"
$ws.Range("E51").Value = "  +0.21%  "
